# Regenerate orders with updated distance/size codes.
# The experiment's Distance codes (D64/D80/D51) and the "far" Size code
# (S30) were renumbered to D69/D86/D55/S31 respectively. These codes are
# embedded as substrings throughout the Condition, Filename_Left,
# Filename_Right, Distance and Size columns (e.g. "Face05_D64_S30" ->
# "Face05_D69_S31", "Fixation_D80_l.png" -> "Fixation_D86_l.png"), so a
# global substring replace over the whole used range reproduces every
# touched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange

$ur.Replace("D64", "D69") | Out-Null
$ur.Replace("D80", "D86") | Out-Null
$ur.Replace("D51", "D55") | Out-Null
$ur.Replace("S30", "S31") | Out-Null
